# [DataTable] ConfigTable structure addition
# - Rename the sheet to match the workbook's purpose
# - Fill in the 4th column's example-value cell (row 4) as part of extending
#   the ConfigTable structure
# - Leave the selection where the author left it when they saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "FarmerSalesTable"

$ws.Range("D4").Value = "z"

$ws.Range("D20").Select()
